# Scheduled market-data refresh for Moogle_Profits workbook.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for
# a set of leve rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
# with freshly pulled marketboard averages.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 859.2
$ws.Range("J103").Value = 924
$ws.Range("L103").Value = 2772
$ws.Range("N103").Value = -3944

$ws.Range("H132").Value = 3940
$ws.Range("I132").Value = 4073.9285
$ws.Range("K132").Value = 12221.7855
$ws.Range("M132").Value = -9691.7855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 712.80646
$ws.Range("I2").Value = 648.6087
$ws.Range("K2").Value = 648.6087
$ws.Range("M2").Value = -535.6087

$ws.Range("H32").Value = 4890.798
$ws.Range("I32").Value = 2275.5657
$ws.Range("J32").Value = 29735.5
$ws.Range("K32").Value = 2275.5657
$ws.Range("L32").Value = 29735.5
$ws.Range("M32").Value = -1988.5657
$ws.Range("N32").Value = -30309.5

$ws.Range("H88").Value = 2615.182
$ws.Range("I88").Value = 1673.75
$ws.Range("J88").Value = 3153.1428
$ws.Range("K88").Value = 1673.75
$ws.Range("L88").Value = 3153.1428
$ws.Range("M88").Value = -1267.75
$ws.Range("N88").Value = -3965.1428

$ws.Range("H91").Value = 2615.182
$ws.Range("I91").Value = 1673.75
$ws.Range("J91").Value = 3153.1428
$ws.Range("K91").Value = 1673.75
$ws.Range("L91").Value = 3153.1428
$ws.Range("M91").Value = -269.75
$ws.Range("N91").Value = -5961.1428

$ws.Range("H116").Value = 712.80646
$ws.Range("I116").Value = 648.6087
$ws.Range("K116").Value = 648.6087
$ws.Range("M116").Value = 1645.3913

$ws.Range("H122").Value = 3994.32
$ws.Range("I122").Value = 2717.95
$ws.Range("K122").Value = 8153.849999999999
$ws.Range("M122").Value = -5703.849999999999

$ws.Range("H133").Value = 37347
$ws.Range("J133").Value = 37347
$ws.Range("L133").Value = 37347
$ws.Range("N133").Value = -42407

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 712.80646
$ws.Range("I3").Value = 648.6087
$ws.Range("K3").Value = 648.6087
$ws.Range("M3").Value = -534.6087

$ws.Range("H94").Value = 1044.2069
$ws.Range("I94").Value = 774.7308
$ws.Range("J94").Value = 3379.6667
$ws.Range("K94").Value = 774.7308
$ws.Range("L94").Value = 3379.6667
$ws.Range("M94").Value = -323.7308
$ws.Range("N94").Value = -4281.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5369.7144
$ws.Range("I31").Value = 2401.1072
$ws.Range("J31").Value = 11306.929
$ws.Range("K31").Value = 2401.1072
$ws.Range("L31").Value = 11306.929
$ws.Range("M31").Value = -2106.1072
$ws.Range("N31").Value = -11896.929

$ws.Range("H34").Value = 5369.7144
$ws.Range("I34").Value = 2401.1072
$ws.Range("J34").Value = 11306.929
$ws.Range("K34").Value = 2401.1072
$ws.Range("L34").Value = 11306.929
$ws.Range("M34").Value = -2199.1072
$ws.Range("N34").Value = -11710.929

$ws.Range("H122").Value = 2478.1177
$ws.Range("I122").Value = 2268.8462
$ws.Range("K122").Value = 6806.5386
$ws.Range("M122").Value = -4356.5386

$ws.Range("H132").Value = 4808.8125
$ws.Range("I132").Value = 2995.5386
$ws.Range("K132").Value = 8986.6158
$ws.Range("M132").Value = -6456.6158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 68.44444
$ws.Range("I38").Value = 69.14286
$ws.Range("K38").Value = 207.42858
$ws.Range("M38").Value = 139.57142

$ws.Range("H113").Value = 1572.25
$ws.Range("J113").Value = 1572.25
$ws.Range("L113").Value = 4716.75
$ws.Range("N113").Value = -9056.75

$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("M124").ClearContents()
$ws.Range("N124").ClearContents()

$ws.Range("H131").Value = 1013114.7
$ws.Range("I131").Value = 1038.6666
$ws.Range("J131").Value = 2025190.6
$ws.Range("K131").Value = 3115.9998
$ws.Range("L131").Value = 6075571.800000001
$ws.Range("M131").Value = 1924.0002
$ws.Range("N131").Value = -6085651.800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 10791.363
$ws.Range("I122").Value = 7141.2
$ws.Range("K122").Value = 21423.6
$ws.Range("M122").Value = -18973.6

$ws.Range("H132").Value = 5098.355
$ws.Range("I132").Value = 3101.6924
$ws.Range("J132").Value = 6540.3887
$ws.Range("K132").Value = 9305.0772
$ws.Range("L132").Value = 19621.1661
$ws.Range("M132").Value = -6775.0772
$ws.Range("N132").Value = -24681.1661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1355.7059
$ws.Range("I16").Value = 1161.7693
$ws.Range("J16").Value = 1986
$ws.Range("K16").Value = 1161.7693
$ws.Range("L16").Value = 1986
$ws.Range("M16").Value = -991.7692999999999
$ws.Range("N16").Value = -2326

$ws.Range("H40").Value = 5740.25
$ws.Range("I40").Value = 3173
$ws.Range("K40").Value = 3173
$ws.Range("M40").Value = -3037

$ws.Range("H46").Value = 2841.5
$ws.Range("I46").Value = 660.7273
$ws.Range("J46").Value = 5022.273
$ws.Range("K46").Value = 660.7273
$ws.Range("L46").Value = 5022.273
$ws.Range("M46").Value = -472.7273
$ws.Range("N46").Value = -5398.273

$ws.Range("H68").Value = 4126
$ws.Range("I68").Value = 4532
$ws.Range("J68").Value = 2366.6667
$ws.Range("K68").Value = 4532
$ws.Range("L68").Value = 2366.6667
$ws.Range("M68").Value = -3783
$ws.Range("N68").Value = -3864.6667

$ws.Range("H71").Value = 4126
$ws.Range("I71").Value = 4532
$ws.Range("J71").Value = 2366.6667
$ws.Range("K71").Value = 22660
$ws.Range("L71").Value = 11833.3335
$ws.Range("M71").Value = -18916
$ws.Range("N71").Value = -19321.3335

$ws.Range("H82").Value = 1573.6666
$ws.Range("I82").Value = 1077.2
$ws.Range("J82").Value = 2194.25
$ws.Range("K82").Value = 1077.2
$ws.Range("L82").Value = 2194.25
$ws.Range("M82").Value = -716.2
$ws.Range("N82").Value = -2916.25

$ws.Range("H85").Value = 1573.6666
$ws.Range("I85").Value = 1077.2
$ws.Range("J85").Value = 2194.25
$ws.Range("K85").Value = 1077.2
$ws.Range("L85").Value = 2194.25
$ws.Range("M85").Value = 170.8
$ws.Range("N85").Value = -4690.25

$ws.Range("H122").Value = 8275.643
$ws.Range("I122").Value = 6975.9
$ws.Range("J122").Value = 11525
$ws.Range("K122").Value = 20927.7
$ws.Range("L122").Value = 34575
$ws.Range("M122").Value = -18477.7
$ws.Range("N122").Value = -39475

$ws.Range("H133").Value = 94116.5
$ws.Range("J133").Value = 94116.5
$ws.Range("L133").Value = 94116.5
$ws.Range("N133").Value = -99176.5

$ws.Range("H136").Value = 8064.2173
$ws.Range("I136").Value = 3225.6667
$ws.Range("K136").Value = 9677.000100000001
$ws.Range("M136").Value = -7127.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3072.889
$ws.Range("I122").Value = 2973
$ws.Range("K122").Value = 8919
$ws.Range("M122").Value = -6469

